{"js": "// Tighten the \"Compact\" paragraph style (used by the numbered/bulleted\n// \"methods\" steps in the body) to a smaller, single-size font with no\n// first-line indent, matching the style's own compact paragraph spacing.\nconst style = context.document.getStyles().getByNameOrNullObject(\"Compact\");\nstyle.load(\"nameLocal\");\nawait context.sync();\n\nif (style.isNullObject) {\n  throw new Error('Paragraph style \"Compact\" was not found in this document.');\n}\n\n// <w:rPr><w:sz w:val=\"20\"/></w:rPr>  -> 20 half-points == 10pt font.\nstyle.font.size = 10;\n\n// <w:pPr>...<w:ind w:firstLine=\"0\"/></w:pPr> -> remove the inherited\n// first-line indent (BodyText normally indents the first line 0.5in).\nstyle.paragraphFormat.firstLineIndent = 0;\n\nawait context.sync();\n", "ps1": "# Tighten the \"Compact\" paragraph style (used by the numbered/bulleted\n# \"methods\" steps in the body) to a smaller, single-size font with no\n# first-line indent, matching the style's own compact paragraph spacing.\n$d = $word.ActiveDocument\n$style = $d.Styles(\"Compact\")\n\n# <w:rPr><w:sz w:val=\"20\"/></w:rPr>  -> 20 half-points == 10pt font.\n$style.Font.Size = 10\n\n# <w:pPr>...<w:ind w:firstLine=\"0\"/></w:pPr> -> remove the inherited\n# first-line indent (BodyText normally indents the first line 0.5in).\n$style.ParagraphFormat.FirstLineIndent = 0\n"}
